# Improvement on Dice coefficient plot
#
# Slide 1 ("Text Box 5") paragraph 3 holds the long Dice-coefficient caption
# as a single run. The canonical edit splits that run into five runs so the
# two misspelled/flagged words ("thresholded" / "thresholding") sit in their
# own runs, with the surrounding text unchanged. This reproduces that run
# split while leaving every other paragraph's text untouched.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)      # "Text Box 5"
$tr = $sh.TextFrame.TextRange

# Paragraph 3 is the big caption paragraph starting "Dice coefficients ..."
$para = $tr.Paragraphs(3, 1)
$run  = $para.Runs(1, 1)

# Re-split the paragraph's single run into five runs, carrying the same
# character formatting forward (InsertAfter inherits the preceding run's
# rPr, matching sz="1000"/Arial/etc. on every fragment).
$run.Text = "Dice coefficients comparing the "
$run2 = $run.InsertAfter("thresholded")
$run3 = $run2.InsertAfter(" positive and negative T-statistic maps computed using each software package and inference method for each of the three reproduced studies.  Dice is the size of the overlapping region of two images divided by the average size of each region. In this context, a Dice coefficient of 1 would indicate perfect agreement between software on the regions of significant activation, whereas a coefficient of 0 would imply that no voxel was declared significant in both packages after ")
$run4 = $run3.InsertAfter("thresholding")
$run5 = $run4.InsertAfter(" the T-statistic images. ")
